$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Date_1 column (A) needs to move from 2026/01/11 -> 2026/01/12.
# These are plain text cells (not real dates), so we must avoid letting
# Excel's autoconvert turn the string into a date serial number. We do
# this by writing the new text into a scratch cell as a formula (so it
# is stored as a literal string), copying it, and pasting-special just
# the values into the target cells - this keeps the cells as plain text
# without touching their number format / style.
$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

$scratch = $ws.Range("Z1")
$scratch.Formula = '="2026/01/12"'
$scratch.Copy()

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)
}

$scratch.ClearContents()
$excel.CutCopyMode = 0
